$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Locate the "Personnes concernées" list (numId 13) that contains
# the Jean-Michel Sordet / Gérard Jaton / Paolo Mariani paragraphs.
# There are earlier (draft / duplicate) occurrences of this same
# text earlier in the document, so we search from the paragraph
# that contains "Avertir le service informatique de Tamedia" and
# walk backwards to find the right "Jean-Michel Sordet" paragraph.
# ---------------------------------------------------------------

$count = $d.Paragraphs.Count
$jeanMichelIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $ptxt = $d.Paragraphs.Item($i).Range.Text
    if ($ptxt -like "*Jean-Michel Sordet*") {
        $jeanMichelIndex = $i
    }
}

# ---------------------------------------------------------------
# 1) Jean-Michel Sordet paragraph : turn the plain e-mail text run
#    into a real hyperlink.
# ---------------------------------------------------------------
$pJM = $d.Paragraphs.Item($jeanMichelIndex)
$rJM = $pJM.Range.Duplicate
$rJM.Find.Execute("jean-michel.sordet@eerv.ch") | Out-Null
$d.Hyperlinks.Add($rJM, "mailto:jean-michel.sordet@eerv.ch", $null, $null, "jean-michel.sordet@eerv.ch") | Out-Null

# ---------------------------------------------------------------
# 2) Insert a brand new paragraph right after the Jean-Michel Sordet
#    one, for Vincent Volet / Marie-Noëlle Aubert.
# ---------------------------------------------------------------
$pJM2 = $d.Paragraphs.Item($jeanMichelIndex)
$endJM = $pJM2.Range.Duplicate
$endJM.Collapse(0)
$endJM.InsertParagraphAfter()

$newIndex = $jeanMichelIndex + 1
$pNew = $d.Paragraphs.Item($newIndex)
$rNew = $pNew.Range.Duplicate
$rNew.Collapse(1)
$rNew.InsertBefore("Vicent Volet et Marie-Noëlle Aubert (Bonne Nouvelle, OIC) : vincent.volet@eerv.ch, marie-noelle.aubert@eerv.ch ")

$pNew2 = $d.Paragraphs.Item($newIndex)
$rVincent = $pNew2.Range.Duplicate
$rVincent.Find.Execute("vincent.volet@eerv.ch") | Out-Null
$d.Hyperlinks.Add($rVincent, "mailto:vincent.volet@eerv.ch", $null, $null, "vincent.volet@eerv.ch") | Out-Null

$pNew3 = $d.Paragraphs.Item($newIndex)
$rMarie = $pNew3.Range.Duplicate
$rMarie.Find.Execute("marie-noelle.aubert@eerv.ch") | Out-Null
$d.Hyperlinks.Add($rMarie, "mailto:marie-noelle.aubert@eerv.ch", $null, $null, "marie-noelle.aubert@eerv.ch") | Out-Null

# ---------------------------------------------------------------
# 3) Gérard Jaton paragraph (now shifted one down) : split the
#    label text run so an empty "_GoBack" bookmark sits in the
#    middle of it, right before "'informatique".
# ---------------------------------------------------------------
$gerardIndex = $newIndex + 1
$pGerard = $d.Paragraphs.Item($gerardIndex)
$rGerard = $pGerard.Range.Duplicate
$rGerard.Find.Execute("Gérard Jaton (responsable de l") | Out-Null
$rGerard.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rGerard) | Out-Null

# ---------------------------------------------------------------
# 4) Paolo Mariani paragraph : replace the legacy HYPERLINK field
#    (fldChar/instrText) + bookmark with a plain w:hyperlink.
#    Easiest reliable way: insert a brand-new paragraph with the
#    wanted content right after the Gérard Jaton paragraph, then
#    delete the old Paolo Mariani paragraph completely.
# ---------------------------------------------------------------
$pGerard2 = $d.Paragraphs.Item($gerardIndex)
$endGerard = $pGerard2.Range.Duplicate
$endGerard.Collapse(0)
$endGerard.InsertParagraphAfter()

$paoloNewIndex = $gerardIndex + 1
$pPaoloNew = $d.Paragraphs.Item($paoloNewIndex)
$rPaoloNew = $pPaoloNew.Range.Duplicate
$rPaoloNew.Collapse(1)
$rPaoloNew.InsertBefore("Paolo Mariani (responsable du BN à l’EERV) : paolo.mariani@eerv.ch")

$pPaoloNew2 = $d.Paragraphs.Item($paoloNewIndex)
$rPaoloEmail = $pPaoloNew2.Range.Duplicate
$rPaoloEmail.Find.Execute("paolo.mariani@eerv.ch") | Out-Null
$d.Hyperlinks.Add($rPaoloEmail, "mailto:paolo.mariani@eerv.ch", $null, $null, "paolo.mariani@eerv.ch") | Out-Null

# The old Paolo Mariani paragraph (with the legacy field-code
# hyperlink) is now right after the new one; delete it completely,
# paragraph mark included.
$paoloOldIndex = $paoloNewIndex + 1
$pPaoloOld = $d.Paragraphs.Item($paoloOldIndex)
$pPaoloOld.Range.Delete()
